$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '26.421.70'
$ws.Range("E2").Value2 = '  +1.02%  '
$ws.Range("D3").Value2 = '1.672.99'
$ws.Range("E3").Value2 = '  +1.11%  '
$ws.Range("E4").Value2 = '  +0.50%  '
$ws.Range("D5").Value2 = '221.75'
$ws.Range("E5").Value2 = '  +1.80%  '
$ws.Range("D6").Value2 = '0.5337'
$ws.Range("E6").Value2 = '  +0.70%  '
$ws.Range("E7").Value2 = '  +0.48%  '
$ws.Range("D8").Value2 = '0.2665'
$ws.Range("E8").Value2 = '  +1.63%  '
$ws.Range("D9").Value2 = '0.06394'
$ws.Range("E9").Value2 = '  +1.17%  '
$ws.Range("D10").Value2 = '20.95'
$ws.Range("E10").Value2 = '  +2.73%  '
$ws.Range("D11").Value2 = '0.07846'
$ws.Range("E11").Value2 = '  +0.46%  '
$ws.Range("D12").Value2 = '4.539'
$ws.Range("E12").Value2 = '  +0.55%  '
$ws.Range("D13").Value2 = '1.683.07'
$ws.Range("E13").Value2 = '  +1.70%  '
$ws.Range("D14").Value2 = '1.900.96'
$ws.Range("E14").Value2 = '  +0.99%  '
$ws.Range("D15").Value2 = '0.5627'
$ws.Range("E15").Value2 = '  +2.52%  '
$ws.Range("D16").Value2 = "0.0{0}8190" -f [char]0x2085
$ws.Range("E16").Value2 = '  +0.51%  '
$ws.Range("D17").Value2 = '66.24'
$ws.Range("E17").Value2 = '  +1.36%  '
$ws.Range("D18").Value2 = '26.410.83'
$ws.Range("E18").Value2 = '  +1.09%  '
$ws.Range("E19").Value2 = '  +0.51%  '
$ws.Range("D20").Value2 = '4.722'
$ws.Range("E20").Value2 = '  +2.69%  '
$ws.Range("D21").Value2 = '198.16'
$ws.Range("E21").Value2 = '  +3.92%  '
$ws.Range("D22").Value2 = '10.32'
$ws.Range("E22").Value2 = '  +2.41%  '
$ws.Range("D23").Value2 = '6.083'
$ws.Range("E23").Value2 = '  +1.41%  '
$ws.Range("D24").Value2 = '1.010'
$ws.Range("E24").Value2 = '  +0.33%  '
$ws.Range("D25").Value2 = '146.70'
$ws.Range("E25").Value2 = '  +0.95%  '
$ws.Range("D26").Value2 = '0.1230'
$ws.Range("E26").Value2 = '  +0.51%  '
$ws.Range("D27").Value2 = '7.250'
$ws.Range("E27").Value2 = '  +0.68%  '
$ws.Range("D28").Value2 = '16.30'
$ws.Range("E28").Value2 = '  +2.01%  '
$ws.Range("D29").Value2 = '1.505'
$ws.Range("E29").Value2 = '  +1.99%  '
$ws.Range("D30").Value2 = '0.05918'
$ws.Range("E30").Value2 = '  +3.71%  '
$ws.Range("D31").Value2 = '1.289'
$ws.Range("E31").Value2 = '  +1.21%  '
$ws.Range("D32").Value2 = '3.562'
$ws.Range("E32").Value2 = '  +0.45%  '
$ws.Range("D33").Value2 = '3.321'
$ws.Range("E33").Value2 = '  +1.71%  '
$ws.Range("D34").Value2 = '1.618'
$ws.Range("E34").Value2 = '  +1.70%  '
$ws.Range("D35").Value2 = '0.9711'
$ws.Range("E35").Value2 = '  +2.46%  '
$ws.Range("D36").Value2 = '2.840'
$ws.Range("E36").Value2 = '  +1.22%  '
$ws.Range("D37").Value2 = '2.435'
$ws.Range("E37").Value2 = '  +0.54%  '
$ws.Range("D38").Value2 = '0.5840'
$ws.Range("E38").Value2 = '  +2.19%  '
$ws.Range("D39").Value2 = '0.01616'
$ws.Range("E39").Value2 = '  +0.56%  '
$ws.Range("D40").Value2 = '1.080.75'
$ws.Range("E40").Value2 = '  +4.08%  '
$ws.Range("D41").Value2 = '5.905'
$ws.Range("E41").Value2 = '  +1.88%  '
$ws.Range("D42").Value2 = '0.8658'
$ws.Range("E42").Value2 = '  +1.90%  '
$ws.Range("E43").Value2 = '  +0.50%  '
$ws.Range("D44").Value2 = '103.71'
$ws.Range("E44").Value2 = '  -0.35%  '
$ws.Range("D45").Value2 = '1.810.97'
$ws.Range("E45").Value2 = '  +0.88%  '
$ws.Range("D46").Value2 = '58.53'
$ws.Range("E46").Value2 = '  +3.19%  '
$ws.Range("D47").Value2 = "0.0{0}106" -f [char]0x2088
$ws.Range("E47").Value2 = '  +0.56%  '
$ws.Range("E48").Value2 = '  +0.65%  '
$ws.Range("E49").Value2 = '  +1.33%  '
$ws.Range("D50").Value2 = '8.002'
$ws.Range("E50").Value2 = '  +2.16%  '
$ws.Range("D51").Value2 = '0.05166'
$ws.Range("E51").Value2 = '  +0.22%  '
